# "adding limits on investment"
#
# Adds a total-investment check formula to both the PredefinedPlantBuilder
# and VariableRenewableOperator sheets, tweaks a couple of
# "ComissionedYear" values, and updates the active sheet/selection.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("PredefinedPlantBuilder")
$ws3 = $wb.Worksheets.Item("VariableRenewableOperator")

# --- PredefinedPlantBuilder (sheet2) -----------------------------------
# New "limit" style total that sums InstalledPowerInMW for a subset of
# producers (Producer1's column G, Producer4's column F, Producer2's column C).
$ws2.Range("I13").Formula = "=G13+F13+C13"

# Clear Producer2's ComissionedYear and bump Producer3's to 2021.
$ws2.Range("D16").Value = ""
$ws2.Range("E16").Value = 2021

# --- VariableRenewableOperator (sheet3) --------------------------------
# New total across several installed-power columns.
$ws3.Range("M6").Formula = "=C6+D6+G6+H6+I6+J6"

# Clear the PVRooftop ComissionedYear and bump WindOn's to 2021.
$ws3.Range("E11").Value = ""
$ws3.Range("F11").Value = 2021

# --- Selections / active sheet ------------------------------------------
$ws3.Range("L26").Select()

$ws2.Activate()
$ws2.Range("E19").Select()

$wb.Save()
